$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the existing "model name" cell format (bold font + border + center/top
# alignment, currently only on A2) to the new A3:A9 cells BEFORE we touch A2's
# own value, so the style-dedup logic can match against the original sparse
# font definition instead of materialising a brand-new (visually identical
# but structurally different) font entry.
$ws.Range("A2").Copy()
$ws.Range("A3:A9").PasteSpecial(-4122)

# Data rows (sheet row 2 .. row 9), columns A..I, in final order
$data = @(
    @("model_11_9_5", 0.7007846717949489, -1.462225904875495, -0.1312398773988326, -0.20869435708999, 0.3311431705951691, 0.2396810352802277, 1.557208299636841, 0.8596953749656677),
    @("model_11_9_3", 0.7072943650870807, -1.997587712560451, 0.05136682605418941, -0.08115956906394528, 0.3239388763904572, 0.2917948663234711, 1.305841207504272, 0.7689851522445679),
    @("model_11_9_7", 0.7188233910103827, -1.336244422466266, -0.09751214380254325, -0.1688549198649807, 0.3111796379089355, 0.2274175733327866, 1.510780453681946, 0.8313592672348022),
    @("model_11_9_6", 0.722536487268424, -1.314859932218669, -0.07463142814022206, -0.1464514933075842, 0.3070703446865082, 0.2253359407186508, 1.479284048080444, 0.8154245615005493),
    @("model_11_9_4", 0.7309833225761887, -1.648607393725317, 0.1355280920293038, 0.02076456883737154, 0.297722190618515, 0.2578240036964417, 1.189988851547241, 0.6964905261993408),
    @("model_11_9_2", 0.8308574778318611, -0.1803146461532472, 0.677948083577345, 0.6211690073908442, 0.1871909350156784, 0.1148956343531609, 0.4433205723762512, 0.2694471776485443),
    @("model_11_9_0", 0.8652791189718037, 0.6072262660691303, 0.9382502808766907, 0.9153023033942906, 0.1490963250398636, 0.03823386132717133, 0.08500158041715622, 0.06024204194545746),
    @("model_11_9_1", 0.8692609953330096, 0.5367574847639057, 0.889064612573792, 0.8654010506304206, 0.1446895599365234, 0.04509352147579193, 0.152708113193512, 0.09573478251695633)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    for ($col = 2; $col -le 9; $col++) {
        $ws.Cells.Item($row, $col).Value = $rowData[$col - 1]
    }
}
